$wb = $excel.ActiveWorkbook

# Sheet index 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("H32").Value = 1549.75
$ws.Range("I32").Value = 1549.75
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1549.75
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1223.75
$ws.Range("H40").Value = 5675.125
$ws.Range("I40").Value = 4580.2
$ws.Range("J40").Value = 7500
$ws.Range("K40").Value = 4580.2
$ws.Range("L40").Value = 7500
$ws.Range("M40").Value = -4405.2
$ws.Range("N40").Value = -7850
$ws.Range("H43").Value = 890.6667
$ws.Range("I43").Value = 890.6667
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 890.6667
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -821.6667
$ws.Range("H121").Value = 1750
$ws.Range("J121").Value = 1750
$ws.Range("L121").Value = 5250
$ws.Range("N121").Value = -8744
$ws.Range("J137").Value = 2489.5
$ws.Range("L137").Value = 7468.5
$ws.Range("N137").Value = -12568.5
$ws.Range("H138").Value = 7605.75
$ws.Range("J138").Value = 7798.516
$ws.Range("L138").Value = 23395.548
$ws.Range("N138").Value = -33675.548
$ws.Range("N32").ClearContents()
$ws.Range("N43").ClearContents()

# Sheet index 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("H17").Value = 12000
$ws.Range("J17").Value = 12000
$ws.Range("L17").Value = 12000
$ws.Range("N17").Value = -12346
$ws.Range("H32").Value = 13639.333
$ws.Range("I32").Value = 11042.143
$ws.Range("J32").Value = 50000
$ws.Range("K32").Value = 11042.143
$ws.Range("L32").Value = 50000
$ws.Range("M32").Value = -10755.143
$ws.Range("N32").Value = -50574
$ws.Range("H61").Value = 7268.2666
$ws.Range("I61").Value = 6771.154
$ws.Range("K61").Value = 6771.154
$ws.Range("M61").Value = -6559.154
$ws.Range("H102").Value = 3091.0588
$ws.Range("J102").Value = 11132.667
$ws.Range("L102").Value = 11132.667
$ws.Range("N102").Value = -14376.667
$ws.Range("H132").Value = 9201.666999999999
$ws.Range("I132").Value = 8642
$ws.Range("J132").Value = 12000
$ws.Range("K132").Value = 25926
$ws.Range("L132").Value = 36000
$ws.Range("M132").Value = -23396
$ws.Range("N132").Value = -41060
$ws.Range("H136").Value = 7268.2666
$ws.Range("I136").Value = 6771.154
$ws.Range("K136").Value = 20313.462
$ws.Range("M136").Value = -17763.462

# Sheet index 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 2179.75
$ws.Range("I20").Value = 1561.8889
$ws.Range("K20").Value = 1561.8889
$ws.Range("M20").Value = -1314.8889

# Sheet index 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("H31").Value = 3049.5715
$ws.Range("I31").Value = 2999.75
$ws.Range("K31").Value = 2999.75
$ws.Range("M31").Value = -2704.75
$ws.Range("H34").Value = 3049.5715
$ws.Range("I34").Value = 2999.75
$ws.Range("K34").Value = 2999.75
$ws.Range("M34").Value = -2797.75
$ws.Range("H68").Value = 51000
$ws.Range("I68").Value = 35000
$ws.Range("K68").Value = 35000
$ws.Range("M68").Value = -34251
$ws.Range("H71").Value = 51000
$ws.Range("I71").Value = 35000
$ws.Range("K71").Value = 105000
$ws.Range("M71").Value = -101256
$ws.Range("H122").Value = 4625.923
$ws.Range("I122").Value = 5373.1113
$ws.Range("J122").Value = 2944.75
$ws.Range("K122").Value = 16119.3339
$ws.Range("L122").Value = 8834.25
$ws.Range("M122").Value = -13669.3339
$ws.Range("N122").Value = -13734.25
$ws.Range("H134").Value = 2809.25
$ws.Range("I134").Value = 2745.8333
$ws.Range("J134").Value = 2999.5
$ws.Range("K134").Value = 8237.499899999999
$ws.Range("L134").Value = 8998.5
$ws.Range("M134").Value = -5702.499899999999
$ws.Range("N134").Value = -14068.5

# Sheet index 5
$ws = $wb.Worksheets.Item(5)
$ws.Range("H3").Value = 3944.25
$ws.Range("I3").Value = 3944.25
$ws.Range("K3").Value = 11832.75
$ws.Range("M3").Value = -11720.75
$ws.Range("H9").Value = 4923
$ws.Range("I9").Value = 4666.3335
$ws.Range("J9").Value = 5000
$ws.Range("K9").Value = 13999.0005
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = -13775.0005
$ws.Range("N9").Value = -15448
$ws.Range("H17").Value = 847.25
$ws.Range("J17").Value = 1500
$ws.Range("L17").Value = 4500
$ws.Range("N17").Value = -4838
$ws.Range("H26").Value = 695
$ws.Range("J26").Value = 816.6667
$ws.Range("L26").Value = 2450.0001
$ws.Range("N26").Value = -3026.0001
$ws.Range("H60").Value = 1214.5714
$ws.Range("I60").Value = 583.6667
$ws.Range("K60").Value = 1751.0001
$ws.Range("M60").Value = -1500.0001
$ws.Range("H134").Value = 4799.75
$ws.Range("I134").Value = 4799.75
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 14399.25
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -9329.25
$ws.Range("H139").Value = 4518.2856
$ws.Range("I139").Value = 2799.1667
$ws.Range("K139").Value = 8397.500100000001
$ws.Range("M139").Value = -3257.500100000001
$ws.Range("H140").Value = 1432532.2
$ws.Range("I140").Value = 1432532.2
$ws.Range("K140").Value = 4297596.6
$ws.Range("M140").Value = -4292416.6
$ws.Range("N134").ClearContents()

# Sheet index 6
$ws = $wb.Worksheets.Item(6)
$ws.Range("H6").Value = 3071.75
$ws.Range("I6").Value = 2899
$ws.Range("J6").Value = 3129.3333
$ws.Range("K6").Value = 2899
$ws.Range("L6").Value = 3129.3333
$ws.Range("M6").Value = -2786
$ws.Range("N6").Value = -3355.3333
$ws.Range("H16").Value = 3071.75
$ws.Range("I16").Value = 2899
$ws.Range("J16").Value = 3129.3333
$ws.Range("K16").Value = 2899
$ws.Range("L16").Value = 3129.3333
$ws.Range("M16").Value = -2649
$ws.Range("N16").Value = -3629.3333
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("H132").Value = 66670636
$ws.Range("I132").Value = 4956.75
$ws.Range("K132").Value = 14870.25
$ws.Range("M132").Value = -12340.25
$ws.Range("M126").ClearContents()

# Sheet index 7
$ws = $wb.Worksheets.Item(7)
$ws.Range("H46").Value = 1235.4
$ws.Range("J46").Value = 1499.3334
$ws.Range("L46").Value = 1499.3334
$ws.Range("N46").Value = -1875.3334
$ws.Range("H100").Value = 3499.5
$ws.Range("I100").Value = 3499.5
$ws.Range("K100").Value = 3499.5
$ws.Range("M100").Value = -2958.5
$ws.Range("H132").Value = 2334.3333
$ws.Range("I132").Value = 2334.3333
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7002.999899999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4472.999899999999
$ws.Range("N132").ClearContents()

# Sheet index 8
$ws = $wb.Worksheets.Item(8)
$ws.Range("H132").Value = 166668740
$ws.Range("I132").Value = 2484.2
$ws.Range("K132").Value = 7452.599999999999
$ws.Range("M132").Value = -4922.599999999999
$ws.Range("H133").Value = 133248.75
$ws.Range("J133").Value = 133248.75
$ws.Range("L133").Value = 133248.75
$ws.Range("N133").Value = -143368.75
$ws.Range("H136").Value = 7492.923
$ws.Range("I136").Value = 8106.75
$ws.Range("K136").Value = 24320.25
$ws.Range("M136").Value = -21770.25
